$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the date header row (row 2) with 4 more dates, reusing the
# existing date-header style from C2 so no duplicate style is created ---
$ws.Range("C2").Copy($ws.Range("D2:G2"))
$ws.Range("D2").Value = 43989
$ws.Range("E2").Value = 43990
$ws.Range("F2").Value = 43991
$ws.Range("G2").Value = 43992

# --- Row 8: new section header "Investigación" (bold, left aligned, no indent) ---
$ws.Range("A3").Copy($ws.Range("A8"))
$ws.Range("A8").Value = "Investigación"
$ws.Range("A8").HorizontalAlignment = -4131

# --- Row 9: task line with indent style (like A4/A5/A7) ---
$ws.Range("A7").Copy($ws.Range("A9"))
$ws.Range("A9").Value = "Aprendizaje por refuerzo (Deep Q-Learning)"

# Duration marker cells reuse the yellow-fill / right-aligned style from C7
$ws.Range("C7").Copy($ws.Range("D9:E9"))
$ws.Range("D9").Value = "4 h."
$ws.Range("E9").Value = "4 h."

$ws.Range("C7").Copy($ws.Range("F9"))
$ws.Range("F9").Value = "3 h."

# --- Row 10: new section header "Implementación Deep Q-Learning" ---
$ws.Range("A8").Copy($ws.Range("A10"))
$ws.Range("A10").Value = "Implementación Deep Q-Learning"

# --- Row 11: task line "Estructuras básicas (estados y agente)" ---
$ws.Range("A7").Copy($ws.Range("A11"))
$ws.Range("A11").Value = "Estructuras básicas (estados y agente)"

$ws.Range("C7").Copy($ws.Range("G11"))
$ws.Range("G11").Value = "2 h."

# --- Column width & selection tweaks to match the final layout ---
# Target stored width is 44.7109375; this runtime quantizes ColumnWidth to
# 1/6-character steps, so 43.85 (-> stored 44.66667) is the closest
# reachable value to the target.
$ws.Columns.Item(1).ColumnWidth = 43.85
$ws.Range("M14").Select()
